$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.720897
$ws.Range("H2").Value = 2.162691
$ws.Range("I2").Value = 0.0284720950782092
$ws.Range("J2").Value = 0.02847209507820921
$ws.Range("M2").Value = 0.06694666666666667
$ws.Range("O2").Value = 0.5098924310779488
$ws.Range("P2").Value = 0.5098924310779488
$ws.Range("Q2").Value = 0.04826165116
$ws.Range("R2").Value = 0.43435486044
$ws.Range("S2").Value = 0.01451770577731059
$ws.Range("T2").Value = 0.01451770577731059

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.720897
$ws.Range("H3").Value = 2.162691
$ws.Range("I3").Value = 0.0284720950782092
$ws.Range("J3").Value = 0.02847209507820921
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.064349
$ws.Range("N3").Value = 0.193047
$ws.Range("O3").Value = 0.4901075689220513
$ws.Range("P3").Value = 0.4901075689220513
$ws.Range("Q3").Value = 0.04638900105300001
$ws.Range("R3").Value = 0.417501009477
$ws.Range("S3").Value = 0.01395438930089861
$ws.Range("T3").Value = 0.01395438930089861

# Row 4
$ws.Range("I4").Value = 0.9467537483444692
$ws.Range("J4").Value = 0.9467537483444693
$ws.Range("M4").Value = 0.06694666666666667
$ws.Range("O4").Value = 0.5098924310779488
$ws.Range("P4").Value = 0.5098924310779488
$ws.Range("S4").Value = 0.482742570375522
$ws.Range("T4").Value = 0.482742570375522

# Row 5
$ws.Range("I5").Value = 0.9467537483444692
$ws.Range("J5").Value = 0.9467537483444693
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.064349
$ws.Range("N5").Value = 0.193047
$ws.Range("O5").Value = 0.4901075689220513
$ws.Range("P5").Value = 0.4901075689220513
$ws.Range("Q5").Value = 1.542526481042
$ws.Range("R5").Value = 13.882738329378
$ws.Range("S5").Value = 0.4640111779689473
$ws.Range("T5").Value = 0.4640111779689474

# Row 6
$ws.Range("G6").Value = 0.6272673333333333
$ws.Range("H6").Value = 1.881802
$ws.Range("I6").Value = 0.0247741565773216
$ws.Range("J6").Value = 0.0247741565773216
$ws.Range("M6").Value = 0.06694666666666667
$ws.Range("O6").Value = 0.5098924310779488
$ws.Range("P6").Value = 0.5098924310779488
$ws.Range("Q6").Value = 0.04199345707555555
$ws.Range("R6").Value = 0.3779411136800001
$ws.Range("S6").Value = 0.01263215492511626
$ws.Range("T6").Value = 0.01263215492511627

# Row 7
$ws.Range("G7").Value = 0.6272673333333333
$ws.Range("H7").Value = 1.881802
$ws.Range("I7").Value = 0.0247741565773216
$ws.Range("J7").Value = 0.0247741565773216
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.064349
$ws.Range("N7").Value = 0.193047
$ws.Range("O7").Value = 0.4901075689220513
$ws.Range("P7").Value = 0.4901075689220513
$ws.Range("Q7").Value = 0.04036402563266667
$ws.Range("R7").Value = 0.363276230694
$ws.Range("S7").Value = 0.01214200165220534
$ws.Range("T7").Value = 0.01214200165220534
